$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marks")

# Set the new Session 5 (Dynamic Prog.) mark for the student
$ws.Range("F4").Value = 3

# Add the comment for the Session 5 column (F5, merged F5:F12)
$ws.Range("F5").Value = "Only dynamic programming but final results are not very good. Keep on!`n"

# Update the active selection to match the new session column
$ws.Range("F5:F12").Select()

